$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$rows = @(
    @(94,  "Butcher#1432", "/get_messages Butcher#1432", "2021-09-12 16:09:07.919000", "ARMA 3 COOP"),
    @(95,  "Butcher#1432", "/get_all",                    "2021-09-13 08:40:23.295000", "ARMA 3 COOP"),
    @(96,  "Butcher#1432", "/get_all",                    "2021-09-13 08:41:13.335000", "ARMA 3 COOP"),
    @(97,  "Butcher#1432", "/get_all",                    "2021-09-13 08:42:32.656000", "ARMA 3 COOP"),
    @(98,  "Butcher#1432", "/get_all",                    "2021-09-13 08:43:03.417000", "ARMA 3 COOP"),
    @(99,  "Butcher#1432", "/get_all",                    "2021-09-13 08:43:30.786000", "ARMA 3 COOP"),
    @(100, "Butcher#1432", "/get_all",                    "2021-09-13 08:44:01.394000", "ARMA 3 COOP"),
    @(101, "Butcher#1432", "/get_all",                    "2021-09-13 08:44:34.845000", "ARMA 3 COOP"),
    @(102, "Butcher#1432", "/get_all",                    "2021-09-13 08:44:52.933000", "ARMA 3 COOP"),
    @(103, "Butcher#1432", "/get_all",                    "2021-09-13 08:45:26.160000", "ARMA 3 COOP"),
    @(104, "Butcher#1432", "/get_all",                    "2021-09-14 07:55:14.929000", "ARMA 3 COOP")
)

$startRow = 82
$templateRow = $startRow - 1

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
    $ws.Cells.Item($r, 5).Value = $data[4]
}

$endRow = $startRow + $rows.Count - 1
$src = $ws.Range("A" + $templateRow + ":E" + $templateRow)
$dst = $ws.Range("A" + $startRow + ":E" + $endRow)
$src.Copy()
$dst.PasteSpecial(-4122)
$excel.CutCopyMode = 0
